# "Fruta / hortaliza, semanal" weekly refresh:
# A new daily record (Fecha 45173, Volumen 40) is inserted at the top of the
# Ciboulette / Vega Modelo de Temuco block (row 349), pushing the existing
# rows 349:422 down by one row (to 350:423). The row that falls off the
# bottom (old row 422) lands in the newly appended row 423.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing block of records (rows 349-422) down by one row,
# making room at row 349 for the new record and growing the used range
# to row 423.
$src = $ws.Range("A349:R422")
$dst = $ws.Range("A350:R423")
$dst.Value = $src.Value2

# The shifted Fecha column (D) needs to keep the original date number
# format (the new bottom row, 423, otherwise loses it).
$ws.Range("D350:D423").NumberFormat = $ws.Range("D349").NumberFormat

# Write the new record into the now-vacated row 349 (Volumen + Fecha are
# the only columns differing from what was already there; the rest of the
# row keeps the same market/category metadata).
$ws.Range("D349").Value = 45173
$ws.Range("J349").Value = 40
